# Update the cryptos price/volume table (columns D = Price, E = Volume(1h))
# with the latest scraped values, per the automated "Updated cryptos list"
# GitHub Actions commit.
#
# Most rows update both the Price (D) and Volume(1h) (E) columns; two rows
# (HuobiToken / TrustWalletToken) only had their Volume(1h) figure refresh
# this run, so only column E is touched for those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "27.324.79"; E = "  +1.31%  " }
    @{ Row = 3; D = "1.826.43"; E = "  +0.51%  " }
    @{ Row = 4; D = "0.9987"; E = "  -0.40%  " }
    @{ Row = 5; D = "313.34"; E = "  +0.82%  " }
    @{ Row = 6; D = "0.9990"; E = "  -0.24%  " }
    @{ Row = 7; D = "0.4473"; E = "  -0.29%  " }
    @{ Row = 8; D = "0.3776"; E = "  +2.38%  " }
    @{ Row = 9; D = "0.07406"; E = "  +1.55%  " }
    @{ Row = 10; D = "0.8805"; E = "  +3.04%  " }
    @{ Row = 11; D = "20.88"; E = "  +0.90%  " }
    @{ Row = 12; D = "1.817.83"; E = "  +0.30%  " }
    @{ Row = 13; D = "6.722"; E = "  +1.48%  " }
    @{ Row = 14; D = "5.438"; E = "  +2.31%  " }
    @{ Row = 15; D = "93.13"; E = "  +1.47%  " }
    @{ Row = 16; D = "0.07062"; E = "  -0.48%  " }
    @{ Row = 17; D = "0.9992"; E = "  -0.48%  " }
    @{ Row = 18; D = "0.000008819"; E = "  +0.83%  " }
    @{ Row = 19; D = "0.9998"; E = "  -0.15%  " }
    @{ Row = 20; D = "15.06"; E = "  +0.94%  " }
    @{ Row = 21; D = "27.318.14"; E = "  +1.11%  " }
    @{ Row = 22; D = "5.347"; E = "  +3.68%  " }
    @{ Row = 23; D = "10.96"; E = "  +0.61%  " }
    @{ Row = 24; D = "1.957"; E = "  -1.45%  " }
    @{ Row = 25; D = "150.97"; E = "  -0.59%  " }
    @{ Row = 26; D = "2.280"; E = "  +2.82%  " }
    @{ Row = 27; D = "18.60"; E = "  +0.80%  " }
    @{ Row = 28; D = "5.351"; E = "  +2.43%  " }
    @{ Row = 29; D = "117.30"; E = "  +0.88%  " }
    @{ Row = 30; D = "0.08874"; E = "  +0.30%  " }
    @{ Row = 31; D = "0.7909"; E = "  +5.58%  " }
    @{ Row = 32; D = "1.198"; E = "  +1.89%  " }
    @{ Row = 33; D = "4.572"; E = "  +3.07%  " }
    @{ Row = 35; D = "0.9992"; E = "  -0.14%  " }
    @{ Row = 37; D = "0.01972"; E = "  +0.68%  " }
    @{ Row = 38; D = "0.05266"; E = "  +0.95%  " }
    @{ Row = 39; D = "7.298"; E = "  +2.70%  " }
    @{ Row = 40; D = "0.5303"; E = "  +0.17%  " }
    @{ Row = 41; D = "2.352"; E = "  +20.10%  " }
    @{ Row = 42; D = "2.876"; E = "  -0.06%  " }
    @{ Row = 43; D = "0.1704"; E = "  +0.51%  " }
    @{ Row = 44; D = "8.638"; E = "  +2.14%  " }
    @{ Row = 45; D = "0.5056"; E = "  -2.92%  " }
    @{ Row = 46; D = "10.60"; E = "  -0.32%  " }
    @{ Row = 47; D = "105.54"; E = "  +0.27%  " }
    @{ Row = 48; D = "1.686"; E = "  +1.47%  " }
    @{ Row = 49; D = "0.9987"; E = "  -0.17%  " }
    @{ Row = 50; D = "0.06386"; E = "  +0.11%  " }
    @{ Row = 51; D = "66.03"; E = "  +5.06%  " }
)

$eOnlyUpdates = @(
    @{ Row = 34; E = "  -1.17%  " }
    @{ Row = 36; E = "  +1.66%  " }
)

foreach ($u in $updates) {
    # Column D values are free-form price strings (e.g. "27.324.79",
    # "0.9992", "1.000"-style figures) that must stay literal text instead
    # of being reinterpreted as numbers (which would silently drop
    # formatting such as trailing zeros or collapse "a.bbb.cc" groupings).
    # Force the cell to Text before assigning, then drop back to General
    # so no stray number-format style sticks around on the cell.
    $dCell = $ws.Cells.Item($u.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
    $dCell.NumberFormat = "general"

    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

foreach ($u in $eOnlyUpdates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
